# Apply the textual changes described by the diff:
# - update the header date
# - update each multiplication expression in the table

$d = $word.ActiveDocument

$replacements = @(
    @{Old = "2024-07-04 Thursday"; New = "2024-07-05 Friday"},
    @{Old = "118×5="; New = "322×6="},
    @{Old = "887×9="; New = "213×6="},
    @{Old = "499×6="; New = "353×2="},
    @{Old = "112×3="; New = "359×6="},
    @{Old = "866×9="; New = "572×5="},
    @{Old = "859×7="; New = "170×6="},
    @{Old = "345×9="; New = "475×2="},
    @{Old = "228×8="; New = "251×7="},
    @{Old = "169×6="; New = "849×2="},
    @{Old = "911×7="; New = "963×6="},
    @{Old = "884×6="; New = "286×9="},
    @{Old = "807×7="; New = "914×2="},
    @{Old = "901×4="; New = "805×5="},
    @{Old = "774×3="; New = "780×3="},
    @{Old = "192×7="; New = "540×3="},
    @{Old = "484×2="; New = "512×3="},
    @{Old = "634×7="; New = "785×9="},
    @{Old = "378×6="; New = "842×4="},
    @{Old = "988×6="; New = "260×5="},
    @{Old = "688×2="; New = "856×3="},
    @{Old = "933×8="; New = "857×3="},
    @{Old = "673×6="; New = "583×5="},
    @{Old = "643×4="; New = "759×4="},
    @{Old = "311×8="; New = "249×5="},
    @{Old = "495×5="; New = "630×6="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $r.New, 2)
}
